# Splash Screen e Design de Tela
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The phone number in B2 gains an extra digit (999999999 -> 9999999999)
$ws.Range("B2").Value = "9999999999"

# Move/extend the active selection from B5 to A3:B3
$ws.Range("A3:B3").Select()
